$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September transaction ("your relationship") was logged on 2024-09-10.
# It belongs at the top of the September log, so insert a fresh row above the
# current row 36 (shifting all following rows, including the August/Broadband
# tail at the bottom, down by one) and populate the new row's details/date.
$ws.Rows("36:36").Insert()

$ws.Range("R36").Value = "your relationship"
$ws.Range("S36").Value = "2024-09-10 11:02:23"
